$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Add an "Assignment Groups 130" paragraph right after the "Kei Giliam ...
# s5270448" line (the blank paragraph that currently separates the team
# roster from the "Date:/Hour:/Name:/Task:" template block), keeping a
# blank paragraph after it (reusing the pre-existing blank paragraph so
# its formatting is untouched).
$introBlank = $d.Paragraphs(6)
$introBlank.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(6)
$newPara.Range.Text = "Assignment Groups 130"

# --- Change 2 -------------------------------------------------------------
# The "Fork ", "Git" (wrapped in spell-check proofErr markers) and " repo"
# runs become a single "Fork Git repo" run.
$d.Content.Find.Execute("Fork Git repo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fork Git repo", 2)

# --- Change 3 ---------------------------------------------------------------
# The "Date:" entry that follows the "Fork Git repo" task (two paragraphs
# later: the task line, then a blank line, then the date line) loses its
# stray <w:lastRenderedPageBreak/>. Locate that paragraph relative to the
# "Fork Git repo" paragraph (rather than a fixed index) since Change 1
# shifts every later paragraph index down by one.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*Fork Git repo*") {
        $target = $para
        break
    }
}
$datePara = $target.Next().Next()
$datePara.Range.Find.Execute("Date:", $true, $false, $false, $false, $false,
                              $true, 1, $false, "Date:", 2)
